# Rewrite the "Micro results" cell of the review table: replace the single
# run (separated by <w:br/> line breaks) with a sequence of per-line
# paragraphs in blue (0000FF) Times New Roman 10pt text, matching the
# updated micro-results summarization.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the "Micro results" row by its label in column 1 so the script
# does not depend on a hard-coded row index.
$targetCell = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $label = $t.Cell($i, 1).Range.Text
    if ($label -like "Micro results*") {
        $targetCell = $t.Cell($i, 2)
        break
    }
}

if ($targetCell -eq $null) {
    throw "Could not find the 'Micro results' row in the table."
}

# Replace the cell's contents wholesale via a raw OOXML payload so each
# result line becomes its own paragraph with the exact run formatting
# from the target revision (color 0000FF, no explicit bold override).
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/08 – BLC – PERIPHERAL–LEFT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/08 – RESP. CULT AND MICRO – **No clear Result**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Low growth of Neisseria spp. and Strep. viridans.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/08 – RESP. CULT AND MICRO – **No clear Result**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Neisseria spp. and Strep viridans &lt;10^4.  </w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>**Explanation:** Both organisms are reported below the clinical significance threshold (typically ≥10^4 CFU/mL for BAL cultures), suggesting colonization or contamination rather than definitive infection. The result is not definitively positive (insufficient growth) but acknowledges their presence, warranting clinical correlation.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/08 – UC – **No clear Result**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Antibiotic usage cautions provided.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/08 – UC – CATHETER SPECIMEN URINE (CSU) NO SIGNIFICANT GROWTH</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>22/08 – UC – **No clear Result**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Advisory on antibiotic use and clinical context needed.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>22/08 – UC – CATHETER SPECIMEN URINE (CSU)CATHETER SITE SWAB NO SIGNIFICANT GROWTH</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>21/08 – BLC – PERIPHERAL–LEFT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>21/08 – RESP. CULT AND MICRO – BRONCHO–ALVEOLAR LAVAGE NO GROWTH</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>20/08 – SARS CORONAVIRUS–2 PCR – **No clear Result**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: MS2 control CT 23.  </w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>(Reason: The result only reports an internal control (MS2) with a CT value, not pathogen–specific findings.)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>20/08 – CPE SCREEN – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: No CPE isolated.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>20/08 – SARS–CoV–2 RNA – Negative</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>20/08 – MRSA SCREEN – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: No MRSA isolated.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetCell.Range.InsertXML($xml) | Out-Null
